# Add 21 new cards to the database.
# Each "card" is a row with Type (column B), Title (D), Description (E),
# and Image location (I). Existing rows 17-36 already have an idx (column A)
# value; new rows 37-38 are brand-new rows appended after the data with no
# idx value. The columns are written per-row in the same order the values
# were originally typed in (some rows have the image path typed before the
# description), so the shared-string table comes out in the same order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cards = @(
    @{ Row = 17; B = 'action'; D = 'Black Hole'; E = 'At least something takes care of the garbage in the universe'; I = 'images/Black Hole.png'; Order = 'DIE' },
    @{ Row = 18; B = 'action'; D = 'Core Drilling Explosion'; E = 'I guess you mined a little too much'; I = 'images/Core Drilling Explosion.png'; Order = 'DEI' },
    @{ Row = 19; B = 'action'; D = 'Dead Planet'; E = 'Life just disappeared'; I = 'images/Dead Planet.png'; Order = 'DEI' },
    @{ Row = 20; B = 'action'; D = 'Fire Storm'; E = 'Save the trees'; I = 'images/Fire Storm.png'; Order = 'DEI' },
    @{ Row = 21; B = 'action'; D = 'gas hurricane'; E = 'As long as you are not near this thing, you will be fine'; I = 'images/gas hurricane.png'; Order = 'DEI' },
    @{ Row = 22; B = 'action'; D = 'Hard Baked'; E = 'Try not to waste too much water'; I = 'images/Hard Baked.png'; Order = 'DEI' },
    @{ Row = 23; B = 'action'; D = 'Ice Age'; E = 'Everyone there are cold hearted'; I = 'images/Ice Age.png'; Order = 'DEI' },
    @{ Row = 24; B = 'action'; D = 'Lighting Storm'; E = 'Magnetic storms can get nasty'; I = 'images/Lighting Storm.png'; Order = 'DEI' },
    @{ Row = 25; B = 'action'; D = 'Nucular Explosion'; E = 'Nucular power is the best way to get energy'; I = 'images/Nucular Explosion.png'; Order = 'DIE' },
    @{ Row = 26; B = 'action'; D = 'Nucular Fission'; E = 'Your planet became a star, now there''s energy'; I = 'images/Nucular Fission.png'; Order = 'DIE' },
    @{ Row = 27; B = 'action'; D = 'Over Exposure'; E = 'One to many trips close to the star'; I = 'images/Over Exposure.png'; Order = 'DEI' },
    @{ Row = 28; B = 'action'; D = 'Plasma Fire'; E = 'Have fun trying to put this one out'; I = 'images/Plasma Fire.png'; Order = 'DEI' },
    @{ Row = 29; B = 'action'; D = 'Super Iodine Explosion'; E = 'Iodide explosives gone wrong'; I = 'imager/Super Iodine Explosion.png'; Order = 'DEI' },
    @{ Row = 30; B = 'action'; D = 'Super Nova'; E = 'Lots and lots of energy'; I = 'images/Super Nova.png'; Order = 'DEI' },
    @{ Row = 31; B = 'action'; D = 'Super Novo'; E = 'Am I drunk, of is the air on fire?'; I = 'images/Super Novo.png'; Order = 'DEI' },
    @{ Row = 32; B = 'action'; D = 'Time Explosion'; E = 'How strange, that explosion seemed to have happened in the future...'; I = 'images/Time Explosion.png'; Order = 'DIE' },
    @{ Row = 33; B = 'action'; D = 'Total Core Meltdown'; E = 'Too much lava, not enough rock'; I = 'images/Total Core Meltdown.png'; Order = 'DEI' },
    @{ Row = 35; B = 'planet'; D = 'Planet Gregren'; E = 'Known for its high amount of gasses'; I = 'images/Planet Gregren.png'; Order = 'DIE' },
    @{ Row = 36; B = 'planet'; D = 'Planet Lightner'; E = 'Perfect place for a colony'; I = 'images/Planet Lightner.png'; Order = 'DEI' },
    @{ Row = 37; B = 'planet'; D = 'Planet Narges'; E = 'A bit cold, but has a great light show'; I = 'images/Planet Narges.png'; Order = 'DEI' },
    @{ Row = 38; B = 'planet'; D = 'Planet Sistene'; E = 'Named for its heavenly look'; I = 'images/Planet Sistene.png'; Order = 'DEI' }
)

$colMap = @{ 'B' = 2; 'D' = 4; 'E' = 5; 'I' = 9 }

foreach ($card in $cards) {
    $r = $card.Row
    $ws.Cells.Item($r, $colMap['B']).Value = $card.B
    foreach ($col in $card.Order.ToCharArray()) {
        $colName = [string]$col
        $ws.Cells.Item($r, $colMap[$colName]).Value = $card[$colName]
    }
}

$ws.Range("I38").Select()
